# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Naranja" (Valencia variety) at
# row 270 of the data table. This pushes the existing rows 270-345 down by
# one (to 271-346) and grows the used range from A1:T345 to A1:T346.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 270, shifting rows 270:345 down
# to 271:346 (and expanding the sheet dimension to A1:T346).
$ws.Rows.Item(270).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A270").Value = 4
$ws.Range("B270").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C270").Value = "Los Lagos"
$ws.Range("D270").Value = 44588
$ws.Range("E270").Value = 10
$ws.Range("F270").Value = "Fruta"
$ws.Range("G270").Value = 100102
$ws.Range("H270").Value = "Cítricos"
$ws.Range("I270").Value = 100102005
$ws.Range("J270").Value = "Naranja"
$ws.Range("K270").Value = "Valencia"
$ws.Range("L270").Value = "Primera"
$ws.Range("M270").Value = 300
$ws.Range("N270").Value = 16000
$ws.Range("O270").Value = 17000
$ws.Range("P270").Value = 16500
$ws.Range("Q270").Value = "$/caja 15 kilos empedrada"
$ws.Range("R270").Value = "Región de O'Higgins"
$ws.Range("S270").Value = 1100
$ws.Range("T270").Value = 15
